# Updated cryptos list — refresh Price (col D) and Volume(1h) (col E)
# values for rows 2-51 to match the latest scrape.
#
# Note: several Price values look numeric (e.g. "1.002", "0.5055") but
# must stay as literal text to match the source data (trailing zeros,
# exact decimal digits, etc.) — prefixing with a leading single-quote
# forces Excel to store them as text instead of auto-converting to a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.055.01'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.874.72'
$ws.Range("E3").Value = '  -2.31%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''319.67'
$ws.Range("E5").Value = '  -3.18%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '''0.5055'
$ws.Range("E7").Value = '  -3.20%  '
$ws.Range("D8").Value = '''0.3965'
$ws.Range("E8").Value = '  -3.03%  '
$ws.Range("D9").Value = '''0.08213'
$ws.Range("E9").Value = '  -3.50%  '
$ws.Range("D10").Value = '''42.11'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("D12").Value = '''23.45'
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("D13").Value = '1.867.40'
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '''6.287'
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '''91.99'
$ws.Range("D18").Value = '''0.00001087'
$ws.Range("E18").Value = '  -2.32%  '
$ws.Range("D19").Value = '''0.06476'
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '30.053.83'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''5.850'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '''11.12'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").Value = '2.081.68'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("D27").Value = '''21.26'
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("D28").Value = '''160.66'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").Value = '''2.233'
$ws.Range("E29").Value = '  -9.33%  '
$ws.Range("D30").Value = '''127.36'
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D31").Value = '''1.088'
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = '''0.1036'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").Value = '''5.951'
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("D34").Value = '''3.690'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").Value = '''0.02438'
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("D36").Value = '''5.237'
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").Value = '''0.06378'
$ws.Range("E37").Value = '  -3.92%  '
$ws.Range("D38").Value = '''0.2139'
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("D39").Value = '''1.174'
$ws.Range("E39").Value = '  -4.94%  '
$ws.Range("D40").Value = '''8.498'
$ws.Range("E40").Value = '  -4.93%  '
$ws.Range("D41").Value = '''0.6295'
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("D42").Value = '''1.214'
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("D43").Value = '''11.31'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("D44").Value = '''13.32'
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("D45").Value = '''0.5914'
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("D46").Value = '''2.107'
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("D47").Value = '''3.625'
$ws.Range("E47").Value = '  -3.74%  '
$ws.Range("D48").Value = '''122.52'
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("E49").Value = '  -3.30%  '
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("E51").Value = '  -4.85%  '
